# TIMRLINKN-28 Only try to add project times for existing users
#
# The underlying XML diff for this commit is mostly incidental churn from
# re-saving the workbook in Excel (fileVersion/revisionPtr/window position),
# plus a few real content changes on the "timr" worksheet:
#   - column widths for D, I, J were adjusted (and I:J split apart)
#   - J5 (Startzeit) and M5 (Endzeit) now store full date+time serials
#     instead of a bare time-of-day fraction
#   - M5 got a dedicated "date + time" number format (a new cellXfs entry)
#   - the active selection on the sheet moved from B3 to M5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width adjustments ---------------------------------------------
# Column D (Aufgabe) narrower: 32 -> 19
$ws.Columns(4).ColumnWidth = 18.166666666666668
# Column I (Notiz) : 10 -> 13.5
$ws.Columns(9).ColumnWidth = 12.666666666666666
# Column J (Startdatum) : 10 -> 23.5 (now distinct from column I)
$ws.Columns(10).ColumnWidth = 22.666666666666668

# --- cell value / formatting updates on row 5 ------------------------------
# J5 "Startzeit": was a bare time fraction (0.675), now the full date+time
# serial (date 43585 + time fraction) while keeping its existing [h]:mm style.
$ws.Range("J5").Value = 43585.675000000003

# M5 "Endzeit": was a bare time fraction under the [h]:mm style, now a full
# date+time serial using a new "date + time" number format.
$ws.Range("M5").NumberFormat = "m/d/yy h:mm"
$ws.Range("M5").Value = 43585.716666666667

# --- selection moved from B3 to M5 -----------------------------------------
$ws.Range("M5").Select()
